$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary of devices")

# Insert a new row at row 12 (existing rows 12-22 shift down to 13-23),
# matching the format of the row that gets pushed down to row 13.
$ws.Rows.Item(12).Insert()
$ws.Range("A13:F13").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)

# Fill in the new "Program counter" row's data
$ws.Cells.Item(12, 1).Value = "Program counter"
$ws.Cells.Item(12, 2).Value = 200
$ws.Cells.Item(12, 3).Value = 0.05
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Formula = "=B12*D12"
$ws.Cells.Item(12, 6).Formula = "=C12*D12"

# Make "summary of devices" the active sheet with E12 selected
$ws.Activate()
$ws.Range("E12").Select()
